$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

# Row 3: LC13004 / Bryan / Lobos
$ws.Range("A3").Value = "LC13004"
$ws.Range("C3").Value = "Lobos"
$ws.Range("B3").Value = "Bryan"
Set-TextValue "F3" "7.76"
Set-TextValue "H3" "7.52"
$ws.Range("D3").Value = 32
$ws.Range("E3").Value = 0
$ws.Range("G3").Value = 2013
$ws.Range("I3").Value = 1

# carnets for rows 4-6 entered as a block
$ws.Range("A4").Value = "SR11038"
$ws.Range("A5").Value = "AM11098"
$ws.Range("A6").Value = "BV13003"

# names/surnames for rows 4-6
$ws.Range("B4").Value = "Rodrigo"
$ws.Range("C4").Value = "Segovia"
$ws.Range("C5").Value = "Motto"
$ws.Range("B5").Value = "Dario"
$ws.Range("B6").Value = "Elias"
$ws.Range("C6").Value = "Barrera"

# Row 7
$ws.Range("A7").Value = "LL13002"
$ws.Range("B7").Value = "Alam"
$ws.Range("C7").Value = "Lopez"

# cum (F) column for rows 4-7
Set-TextValue "F4" "7.5"
Set-TextValue "F5" "7.6"
Set-TextValue "F6" "8.3"
Set-TextValue "F7" "7.5"

# remaining numeric columns
$ws.Range("D4").Value = 34
$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 2011
$ws.Range("H4").Value = 8
$ws.Range("I4").Value = 1

$ws.Range("D5").Value = 33
$ws.Range("E5").Value = 2
$ws.Range("G5").Value = 2011
$ws.Range("H5").Value = 7
$ws.Range("I5").Value = 1

$ws.Range("D6").Value = 32
$ws.Range("E6").Value = 0
$ws.Range("G6").Value = 2013
$ws.Range("H6").Value = 9
$ws.Range("I6").Value = 1

$ws.Range("D7").Value = 32
$ws.Range("E7").Value = 0
$ws.Range("G7").Value = 2013
$ws.Range("H7").Value = 7
$ws.Range("I7").Value = 1

# Match the author's final selection/active cell
$ws.Range("H3").Select() | Out-Null
